$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B=0.2707932045371706; C=0.07041530944768226; D=0.07796443285583621; E=0.1450545770784402; G=1.415231341819094; H=1.258408207658192; I=1.100024396642105; K=0.3325312292098772; M=0.2439005126269436 }
    3  = @{ B=0.2449364021153713; C=0.06252644056040424; D=0.07079448176847336; E=0.1339653914798689; G=1.382884738630935; H=1.247665512840967; I=1.086065402385685;  K=0.2991984774104708; M=0.2224837851026891 }
    4  = @{ B=0.229201767033004;  C=0.05770978266350824; D=0.06642964103181725; E=0.1272419048668993; G=1.363664709184462; H=1.241537666044167; I=1.077912021176815;  K=0.2788943508040092; M=0.2094705641491146 }
    5  = @{ B=0.2228253649037129; C=0.05575366887325117; D=0.06466030455662519; E=0.1245231929367989; G=1.355992870654617; H=1.23915804690418;  I=1.07469411306942;   K=0.2706609335363197; M=0.2042016281893382 }
    6  = @{ B=0.2217687177334255; C=0.05542925992720882; D=0.06436707164199618; E=0.1240730226592461; G=1.354728640360278; H=1.238770005361346; I=1.074166093755117;  K=0.2692962339492055; M=0.2033287729081223 }
    7  = @{ B=0.2291156285444345; C=0.05768337479177887; D=0.06640574128350352; E=0.1272051540581245; G=1.363560595222623; H=1.241505098085412; I=1.077868200008211;  K=0.2787831475264255; M=0.2093993679666326 }
    8  = @{ B=0.2618483891038466; C=0.06768952509757753; D=0.0754843994926091;  E=0.1412131393834457; G=1.403944786646463; H=1.254606818520216; I=1.095124483637719;  K=0.3210043125263553; M=0.2364874413479967 }
    9  = @{ B=0.3271643263728379; C=0.08753294229271091; D=0.09358949171075892; E=0.169374300212489;  G=1.488259540268871; H=1.284026102124074; I=1.132293791603672;  K=0.4050987536153059; M=0.2907096796980611 }
    10 = @{ B=0.3758500981325028; C=0.1022567331665982;  D=0.1070822938558393;  E=0.1905082112568053; G=1.553385324618034; H=1.307932319162489; I=1.161659472758245;  K=0.4676967949066579; M=0.331247654771289  }
    11 = @{ B=0.3981528974620687; C=0.1089886493913923;  D=0.1132635798590513;  E=0.200223928008036;  G=1.583716066091966; H=1.319310070260912; I=1.175471585339366;  K=0.4963562299581099; M=0.3498481128838193 }
    12 = @{ B=0.406620849029764;  C=0.111542889243367;   D=0.1156106032799187;  E=0.2039180146194184; G=1.595303749058445; H=1.323691115232862; I=1.180767510698885;  K=0.5072354535693648; M=0.3569150025379955 }
    13 = @{ B=0.4047961290739863; C=0.1109925632479474;  D=0.1151048482425523;  E=0.2031217567342125; G=1.592803581677288; H=1.32274434865684;  I=1.179624016012369;  K=0.5048912360694828; M=0.355391979367667  }
    14 = @{ B=0.3988491125674898; C=0.1091986871450672;  D=0.1134565441208366;  E=0.2005275411763066; G=1.584667341380168; H=1.319669046260515; I=1.17590596793849;   K=0.4972507378237481; M=0.3504290411090949 }
    15 = @{ B=0.3952093065504982; C=0.1081005422269641;  D=0.1124477334389979;  E=0.1989404663666079; G=1.579696976804854; H=1.317794788766776; I=1.173637107141161;  K=0.4925741703525262; M=0.3473921432018088 }
    16 = @{ B=0.3743956862980156; C=0.1018174809291281;  D=0.1066792105784486;  E=0.1898753417721508; G=1.551417384061295; H=1.307198893228019; I=1.160765975617466;  K=0.46582753514096;   M=0.3300353118825896 }
    17 = @{ B=0.3616670178137724; C=0.09797181795619281; D=0.1031515607142239;  E=0.1843404786012002; G=1.534249795322467; H=1.300827583992543; I=1.152986367534211;  K=0.4494663953275335; M=0.3194285720002412 }
    18 = @{ B=0.3543604496476007; C=0.0957630862573069;  D=0.101126625092391;   E=0.1811665339438804; G=1.524441736224048; H=1.297210276090226; I=1.148554408058587;  K=0.4400731407409921; M=0.3133428497120079 }
    19 = @{ B=0.3518890826020709; C=0.09501579219175937; D=0.1004417144345666;  E=0.1800935235086882; G=1.521132257324325; H=1.295993634130809; I=1.147061139523927;  K=0.4368957022994664; M=0.3112848930481178 }
    20 = @{ B=0.3630204925163696; C=0.09838086411895119; D=0.1035266627709461;  E=0.1849286822889624; G=1.536070448856691; H=1.301500923477704; I=1.153810101764549;  K=0.4512062811251951; M=0.3205561231625396 }
    21 = @{ B=0.4005952883003943; C=0.1097254551139315;  D=0.1139405193156193;  E=0.2012891167366533; G=1.587054375664763; H=1.320570366458725; I=1.176996266112141;  K=0.4994942163564247; M=0.3518861411398646 }
    22 = @{ B=0.4252829652259891; C=0.117169056628768;    D=0.1207833582590894;  E=0.2120689281097654; G=1.620970727830894; H=1.333456272576655; I=1.192532173928271;  K=0.5312078182370783; M=0.3724980984433515 }
    23 = @{ B=0.4120947503180901; C=0.1131935489414388;  D=0.1171278168573622;  E=0.2063074442938131; G=1.602814201887838; H=1.326540040929927; I=1.184205256691556;  K=0.5142674664892581; M=0.3614845467953671 }
    24 = @{ B=0.3624085513684463; C=0.09819592753459006; D=0.1033570693048489;  E=0.1846627301326009; G=1.535247138888877; H=1.301196364523662; I=1.153437564941314;  K=0.4504196383338979; M=0.320046319209851  }
    25 = @{ B=0.3093728478358457; C=0.08214005672331837; D=0.08865854986626687; E=0.1616796273430907; G=1.464895941172131; H=1.275666375436458; I=1.121879199661748;  K=0.3822076719683025; M=0.2364874413479967 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
